# Scheduled runner update: refresh cached market-board price/profit figures
# (columns H:N) across several rows on each profession sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 914.125
$ws.Range("J33").Value = 499
$ws.Range("L33").Value = 499
$ws.Range("N33").Value = -957
$ws.Range("H58").Value = 5565.5557
$ws.Range("I58").Value = 1364
$ws.Range("J58").Value = 7666.3335
$ws.Range("K58").Value = 4092
$ws.Range("L58").Value = 22999.0005
$ws.Range("M58").Value = -3942
$ws.Range("N58").Value = -23299.0005
$ws.Range("H74").Value = 50006540
$ws.Range("I74").Value = 166669330
$ws.Range("K74").Value = 166669330
$ws.Range("M74").Value = -166668394
$ws.Range("H77").Value = 50006540
$ws.Range("I77").Value = 166669330
$ws.Range("K77").Value = 833346650
$ws.Range("M77").Value = -833341970
$ws.Range("H121").Value = 6000
$ws.Range("J121").Value = 6000
$ws.Range("L121").Value = 18000
$ws.Range("N121").Value = -21494
$ws.Range("H132").Value = 2625.7742
$ws.Range("I132").Value = 2625.7742
$ws.Range("K132").Value = 7877.3226
$ws.Range("M132").Value = -5347.3226
$ws.Range("H135").Value = 770674.75
$ws.Range("I135").Value = 1111946.9
$ws.Range("K135").Value = 10007522.1
$ws.Range("M135").Value = -10004987.1
$ws.Range("H137").Value = 3567.5625
$ws.Range("I137").Value = 7696.6665
$ws.Range("J137").Value = 2614.6924
$ws.Range("K137").Value = 23089.9995
$ws.Range("L137").Value = 7844.0772
$ws.Range("M137").Value = -20539.9995
$ws.Range("N137").Value = -12944.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7478.8213
$ws.Range("I61").Value = 2898.4707
$ws.Range("J61").Value = 14557.546
$ws.Range("K61").Value = 2898.4707
$ws.Range("L61").Value = 14557.546
$ws.Range("M61").Value = -2686.4707
$ws.Range("N61").Value = -14981.546
$ws.Range("H132").Value = 7484.3237
$ws.Range("I132").Value = 6118.952
$ws.Range("J132").Value = 9689.923000000001
$ws.Range("K132").Value = 18356.856
$ws.Range("L132").Value = 29069.769
$ws.Range("M132").Value = -15826.856
$ws.Range("N132").Value = -34129.769
$ws.Range("H136").Value = 7478.8213
$ws.Range("I136").Value = 2898.4707
$ws.Range("J136").Value = 14557.546
$ws.Range("K136").Value = 8695.4121
$ws.Range("L136").Value = 43672.638
$ws.Range("M136").Value = -6145.4121
$ws.Range("N136").Value = -48772.638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 45810.668
$ws.Range("I50").Value = 39500
$ws.Range("J50").Value = 48966
$ws.Range("K50").Value = 39500
$ws.Range("L50").Value = 48966
$ws.Range("M50").Value = -38926
$ws.Range("N50").Value = -50114
$ws.Range("H105").Value = 5034
$ws.Range("I105").Value = 2450.25
$ws.Range("K105").Value = 2450.25
$ws.Range("M105").Value = -703.25
$ws.Range("H134").Value = 6977.4546
$ws.Range("I134").Value = 3586.1875
$ws.Range("K134").Value = 10758.5625
$ws.Range("M134").Value = -8223.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2828.5625
$ws.Range("I16").Value = 1217.75
$ws.Range("K16").Value = 1217.75
$ws.Range("M16").Value = -930.75
$ws.Range("H22").Value = 559.8
$ws.Range("J22").Value = 699.5
$ws.Range("L22").Value = 699.5
$ws.Range("N22").Value = -1399.5
$ws.Range("H31").Value = 9767
$ws.Range("I31").Value = 4386.7144
$ws.Range("J31").Value = 14788.6
$ws.Range("K31").Value = 4386.7144
$ws.Range("L31").Value = 14788.6
$ws.Range("M31").Value = -4091.7144
$ws.Range("N31").Value = -15378.6
$ws.Range("H34").Value = 9767
$ws.Range("I34").Value = 4386.7144
$ws.Range("J34").Value = 14788.6
$ws.Range("K34").Value = 4386.7144
$ws.Range("L34").Value = 14788.6
$ws.Range("M34").Value = -4184.7144
$ws.Range("N34").Value = -15192.6
$ws.Range("H113").Value = 2828.5625
$ws.Range("I113").Value = 1217.75
$ws.Range("K113").Value = 1217.75
$ws.Range("M113").Value = 952.25
$ws.Range("H137").Value = 48750
$ws.Range("J137").Value = 48750
$ws.Range("L137").Value = 48750
$ws.Range("N137").Value = -58950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 555
$ws.Range("J7").Value = 606.25
$ws.Range("L7").Value = 1818.75
$ws.Range("N7").Value = -2042.75
$ws.Range("H113").Value = 3217.7693
$ws.Range("I113").Value = 905.4286
$ws.Range("K113").Value = 2716.2858
$ws.Range("M113").Value = -546.2857999999997
$ws.Range("H129").Value = 11178707
$ws.Range("J129").Value = 23953274
$ws.Range("L129").Value = 71859822
$ws.Range("N129").Value = -71869822

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11999.833
$ws.Range("J70").Value = 11999.833
$ws.Range("L70").Value = 11999.833
$ws.Range("N70").Value = -12539.833
$ws.Range("H73").Value = 11999.833
$ws.Range("J73").Value = 11999.833
$ws.Range("L73").Value = 11999.833
$ws.Range("N73").Value = -13871.833
$ws.Range("H113").Value = 4516.9443
$ws.Range("I113").Value = 2180.6
$ws.Range("K113").Value = 2180.6
$ws.Range("M113").Value = -10.59999999999991

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1232.45
$ws.Range("I16").Value = 1232.45
$ws.Range("K16").Value = 1232.45
$ws.Range("M16").Value = -1062.45
$ws.Range("H68").Value = 90911440
$ws.Range("I68").Value = 142858430
$ws.Range("J68").Value = 4216.5
$ws.Range("K68").Value = 142858430
$ws.Range("L68").Value = 4216.5
$ws.Range("M68").Value = -142857681
$ws.Range("N68").Value = -5714.5
$ws.Range("H71").Value = 90911440
$ws.Range("I71").Value = 142858430
$ws.Range("J71").Value = 4216.5
$ws.Range("K71").Value = 714292150
$ws.Range("L71").Value = 21082.5
$ws.Range("M71").Value = -714288406
$ws.Range("N71").Value = -28570.5
$ws.Range("H122").Value = 3405.9443
$ws.Range("J122").Value = 4665.4165
$ws.Range("L122").Value = 13996.2495
$ws.Range("N122").Value = -18896.2495
$ws.Range("H132").Value = 6400.6978
$ws.Range("I132").Value = 3438.6155
$ws.Range("K132").Value = 10315.8465
$ws.Range("M132").Value = -7785.8465
$ws.Range("H136").Value = 9651.852999999999
$ws.Range("J136").Value = 11475.16
$ws.Range("L136").Value = 34425.48
$ws.Range("N136").Value = -39525.48
$ws.Range("H141").Value = 29999
$ws.Range("J141").Value = 29999
$ws.Range("L141").Value = 29999
$ws.Range("N141").Value = -40359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8191.3335
$ws.Range("I62").Value = 7337.1665
$ws.Range("J62").Value = 9899.666999999999
$ws.Range("K62").Value = 7337.1665
$ws.Range("L62").Value = 9899.666999999999
$ws.Range("M62").Value = -6713.1665
$ws.Range("N62").Value = -11147.667
$ws.Range("H65").Value = 8191.3335
$ws.Range("I65").Value = 7337.1665
$ws.Range("J65").Value = 9899.666999999999
$ws.Range("K65").Value = 36685.8325
$ws.Range("L65").Value = 49498.335
$ws.Range("M65").Value = -33565.8325
$ws.Range("N65").Value = -55738.335
$ws.Range("H107").Value = 1229.6
$ws.Range("I107").Value = 1259.4
$ws.Range("J107").Value = 1199.8
$ws.Range("K107").Value = 3778.2
$ws.Range("L107").Value = 3599.4
$ws.Range("M107").Value = -1858.2
$ws.Range("N107").Value = -7439.4
$ws.Range("H126").Value = 37041600
$ws.Range("I126").Value = 55559060
$ws.Range("K126").Value = 166677180
$ws.Range("M126").Value = -166674710
